# Course schedule: mark homework as assigned ("YES") for the Lists wrap-up
# lesson, the Array/List Algorithms lessons, and the newly added
# Dictionaries lab + exercises rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G17").Value = "YES"
$ws.Range("G18").Value = "YES"
$ws.Range("G19").Value = "YES"
$ws.Range("G21").Value = "YES"
$ws.Range("G22").Value = "YES"

# Rows 18 and 19 pick up the alternating shaded-row look already used by
# the neighboring lesson rows (15 and 13 respectively) - copy their
# formatting across without touching the underlying formulas/values.
$ws.Range("A15:G15").Copy()
$ws.Range("A18:G18").PasteSpecial(-4122)

$ws.Range("A13:G13").Copy()
$ws.Range("A19:G19").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Leave the view where the edit ended up: Dictionaries rows selected.
$ws.Range("G21:G22").Select()
